# Add the new rows (111-138) of "havana club", "brugal" and "puntacana"
# products to the "precios" sheet, continuing directly after the existing
# data (which currently ends at row 110).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("precios")

# Columns, in order: A=pais B=tipo C=tienda D=empresa E=marca F=sku
#                     G=precio H=ml I=grado
$newRows = @(
    ,@('francia','ron','la maison du whisky','havana club','havana club','havana club 7',33.5,700,40)
    ,@('francia','ron','la maison du whisky','havana club','havana club','havana club seleccion',52.5,700,45)
    ,@('francia','ron','la maison du whisky','havana club','havana club','havana club 3',20.9,700,40)
    ,@('francia','ron','rhum attitude','havana club','havana club','havana club anejo especial',22.9,700,40)
    ,@('francia','ron','rhum attitude','havana club','havana club','havana club 3',20.9,700,40)
    ,@('francia','ron','la maison du whisky','havana club','havana club','havana club seleccion',54.9,700,45)
    ,@('francia','ron','urban drinks','havana club','havana club','havana club 3',21.9,700,40)
    ,@('francia','ron','urban drinks','havana club','havana club','havana club 7',31.9,700,40)
    ,@('francia','ron','urban drinks','havana club','havana club','havana club seleccion',52.9,700,45)
    ,@('francia','ron','urban drinks','havana club','havana club','havana club anejo especial',27.9,1000,40)
    ,@('francia','ron','la maison du whisky','brugal','brugal','brugal 1888',49.9,700,40)
    ,@('francia','ron','la maison du whisky','brugal','brugal','brugal anejo',23.9,700,38)
    ,@('francia','ron','la maison du whisky','brugal','brugal','brugal blanco',24,700,40)
    ,@('francia','ron','urban drinks','brugal','brugal','brugal anejo',18.9,700,38)
    ,@('francia','ron','urban drinks','brugal','brugal','brugal 1888',42.9,700,40)
    ,@('francia','ron','urban drinks','brugal','brugal','brugal blanco',20.9,700,40)
    ,@('francia','ron','urban drinks','brugal','brugal','brugal extra viejo',31.9,700,38)
    ,@('francia','ron','rhum attitude','brugal','brugal','brugal anejo',22.9,700,38)
    ,@('francia','ron','rhum attitude','brugal','brugal','brugal 1888',49.9,700,40)
    ,@('francia','ron','rhum attitude','brugal','brugal','brugal extra viejo',31.9,700,38)
    ,@('francia','ron','licorea','brugal','brugal','brugal leyenda 5',98.15,700,38)
    ,@('francia','ron','excellence rhum','brugal','brugal','brugal leyenda',71,700,38)
    ,@('francia','ron','excellence rhum','brugal','brugal','brugal anejo',25.99,700,38)
    ,@('francia','ron','excellence rhum','brugal','brugal','brugal blanco',25,700,40)
    ,@('francia','ron','excellence rhum','brugal','brugal','brugal 1888',45,700,40)
    ,@('francia','ron','excellence rhum','brugal','brugal','brugal reserva xv',39,700,40)
    ,@('francia','ron','rhum attitude','oliver','puntacana','puntacana xox',88.9,700,40)
    ,@('francia','ron','rhum attitude','oliver','puntacana','puntacana tesoro',54.9,700,40)
)

$startRow = 111
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    if ($r -eq 131) {
        # This row's "sku" (F, brugal leyenda 5) was typed in before its
        # "tienda" (C, licorea), so the shared-string table order differs
        # from a plain left-to-right fill for this one row.
        $colOrder = @(1, 2, 6, 3, 4, 5, 7, 8, 9)
    } else {
        $colOrder = @(1, 2, 3, 4, 5, 6, 7, 8, 9)
    }
    foreach ($c in $colOrder) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# Restore the cosmetic view state Excel records on save: active cell at the
# bottom of the newly-entered data, frozen header row, and the pane scrolled
# so the new rows are visible.
$aw = $excel.ActiveWindow
$aw.FreezePanes = $false
$ws.Range("A2").Select()
$aw.FreezePanes = $true
$ws.Range("I138").Select()
